# "computed for year and car age"
#
# 1. Fix three mistyped Car IDs in column A (rows 7, 15, 38).
# 2. Add two computed columns:
#       F = 2-digit "model year" code parsed out of the Car ID (MID(A,3,2))
#       G = car age derived from that 2-digit year code
#    for every data row (2-53).
# 3. Leave the selection on C38 (matches the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("car inventory")

# --- Fix typos in the Car ID column -----------------------------------
$ws.Range("A7").Value = "FD06FCS006"
$ws.Range("A15").Value = "GM09CMR014"
$ws.Range("A38").Value = "HO05ODY037"

# --- New computed columns F (year code) and G (age) --------------------
# Row 2 is entered individually first (its own, non-shared formula),
# then rows 3:53 are filled as a block, mirroring how the neighboring
# B:E helper columns were originally built.
$ws.Range("F2").Formula = "=MID(A2,3,2)"
$ws.Range("G2").Formula = "=IF(25-F2<0,100-F2+25,25-F2)"

$ws.Range("F3:F53").Formula = "=MID(A3,3,2)"
$ws.Range("G3:G53").Formula = "=IF(25-F3<0,100-F3+25,25-F3)"

# --- Restore the selection shown in the saved workbook ------------------
$ws.Range("C38").Select() | Out-Null
